# Added activation costs to optimization cost calculation
# The fuel_price column (N) values for rows 3-26 are multiplied by 3.6
# to account for activation costs in the optimization cost calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 14)  # Column N is the 14th column
    $current = $cell.Value2
    $cell.Value2 = $current * 3.6
}
